$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overview" (sheet1 / table "Overview")
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$rowOverview = $loOverview.ListRows.Add()

$wsOverview.Range("A4").Value = "46febb8d-3c24-4e9a-9f73-67dbb54e4a65.md"
$wsOverview.Range("B4").Value = "e2e\46febb8d-3c24-4e9a-9f73-67dbb54e4a65.md"
$wsOverview.Range("C4").Value = ".md"
$wsOverview.Range("E4").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F4").Value = "Handed back: in sync with en-US"
$wsOverview.Range("G4").Value = "2016-12-16 09:40:22"
$wsOverview.Range("G4").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$hlOverview = $wsOverview.Hyperlinks.Add(
    $wsOverview.Range("B4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/10e063f400618e79060f9a24fe61fc457aa184b7/e2e/46febb8d-3c24-4e9a-9f73-67dbb54e4a65.md",
    "",
    "",
    "e2e\46febb8d-3c24-4e9a-9f73-67dbb54e4a65.md"
)

# ---------------------------------------------------------------------------
# Sheet "zh-cn" (sheet2 / table "zh-cn")
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$loZhCn = $wsZhCn.ListObjects.Item(1)
$rowZhCn = $loZhCn.ListRows.Add()

$wsZhCn.Range("A4").Value = "46febb8d-3c24-4e9a-9f73-67dbb54e4a65.md"
$wsZhCn.Range("B4").Value = ".md"
$wsZhCn.Range("C4").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("D4").Value = "e2e"
$wsZhCn.Range("E4").Value = "ht"
$wsZhCn.Range("F4").Value = "'True"
$wsZhCn.Range("G4").Value = "46febb8d-3c24-4e9a-9f73-67dbb54e4a65.6478912206ffe8134a4b6eb141c667d73b845d25.zh-cn.xlf"
$wsZhCn.Range("H4").Value = "2016-12-16 09:40:08"
$wsZhCn.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("J4").Value = "46febb8d-3c24-4e9a-9f73-67dbb54e4a65.md"
$wsZhCn.Range("K4").Value = "46febb8d-3c24-4e9a-9f73-67dbb54e4a65.6478912206ffe8134a4b6eb141c667d73b845d25.zh-cn.xlf"
$wsZhCn.Range("L4").Value = "2016-12-16 09:41:02"
$wsZhCn.Range("L4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("O4").Value = "'True"
$wsZhCn.Range("Q4").Value = "'False"

$hlZhCnA = $wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("A4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/10e063f400618e79060f9a24fe61fc457aa184b7/e2e/46febb8d-3c24-4e9a-9f73-67dbb54e4a65.md",
    "",
    "",
    "46febb8d-3c24-4e9a-9f73-67dbb54e4a65.md"
)
$hlZhCnJ = $wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("J4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/6478912206ffe8134a4b6eb141c667d73b845d25/e2e/46febb8d-3c24-4e9a-9f73-67dbb54e4a65.md",
    "",
    "",
    "46febb8d-3c24-4e9a-9f73-67dbb54e4a65.md"
)

# ---------------------------------------------------------------------------
# Sheet "de-de" (sheet3 / table "de-de")
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$loDeDe = $wsDeDe.ListObjects.Item(1)
$rowDeDe = $loDeDe.ListRows.Add()

$wsDeDe.Range("A4").Value = "46febb8d-3c24-4e9a-9f73-67dbb54e4a65.md"
$wsDeDe.Range("B4").Value = ".md"
$wsDeDe.Range("C4").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("D4").Value = "e2e"
$wsDeDe.Range("E4").Value = "ht"
$wsDeDe.Range("F4").Value = "'True"
$wsDeDe.Range("G4").Value = "46febb8d-3c24-4e9a-9f73-67dbb54e4a65.6478912206ffe8134a4b6eb141c667d73b845d25.de-de.xlf"
$wsDeDe.Range("H4").Value = "2016-12-16 09:40:22"
$wsDeDe.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("J4").Value = "46febb8d-3c24-4e9a-9f73-67dbb54e4a65.md"
$wsDeDe.Range("K4").Value = "46febb8d-3c24-4e9a-9f73-67dbb54e4a65.6478912206ffe8134a4b6eb141c667d73b845d25.de-de.xlf"
$wsDeDe.Range("L4").Value = "2016-12-16 09:41:21"
$wsDeDe.Range("L4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("O4").Value = "'True"
$wsDeDe.Range("Q4").Value = "'False"

$hlDeDeA = $wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("A4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/10e063f400618e79060f9a24fe61fc457aa184b7/e2e/46febb8d-3c24-4e9a-9f73-67dbb54e4a65.md",
    "",
    "",
    "46febb8d-3c24-4e9a-9f73-67dbb54e4a65.md"
)
$hlDeDeJ = $wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("J4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/6478912206ffe8134a4b6eb141c667d73b845d25/e2e/46febb8d-3c24-4e9a-9f73-67dbb54e4a65.md",
    "",
    "",
    "46febb8d-3c24-4e9a-9f73-67dbb54e4a65.md"
)
